# "modifictation titre graphique et position graphique"
#  - Fix the chart title text for the "part_tps_partiel" row (column C,
#    row 6): it previously held an unrelated label; correct it to read
#    "Proportion de jeunes à temps partiel".
#  - Move the active selection on the sheet to D11 (reflecting where the
#    user last clicked / the new "position graphique").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "Proportion de jeunes à temps partiel"

$ws.Range("D11").Select() | Out-Null
